$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.601.26"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "1.858.80"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6941"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07696"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3060"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07772"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.148"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "1.852.86"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "91.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6924"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.565"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "29.480.51"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008288"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "2.102.81"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.614"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1499"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.928"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.533"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.254"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.178"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05160"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7706"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.893"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.153"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.688"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "1.334.86"
$ws.Range("E38").Value = "  +8.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01873"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.729"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9699"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.813"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000126"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.795"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").Value = "2.001.18"
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5217"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.779"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.964"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.93%  "
